# Apply the CDL_limit.xlsx changes:
#  1. Remove the CDLTotalFuelConsumed and CDLTotalLoadCycles parameter rows.
#  2. Widen CDLCurrentGear max_value from 12 to 17.
#  3. Widen the four strut-pressure sensors' range from 500-1500 to 0-30000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two retired CDL rows (by name, so row numbers don't matter) ---
$namesToDelete = @("CDLTotalFuelConsumed", "CDLTotalLoadCycles")
foreach ($name in $namesToDelete) {
    $cell = $ws.Columns.Item(1).Find($name)
    if ($cell) {
        $cell.EntireRow.Delete()
    }
}

# --- 2. CDLCurrentGear: max_value 12 -> 17 ---
$cell = $ws.Columns.Item(1).Find("CDLCurrentGear")
if ($cell) {
    $r = $cell.Row
    $ws.Cells.Item($r, 3).Value = 17
}

# --- 3. Strut pressure sensors: min 500 -> 0, max 1500 -> 30000 ---
$strutSensors = @(
    "CDLLeftFrontStrutPressure",
    "CDLLeftRearStrutPressure",
    "CDLRightFrontStrutPressure",
    "CDLRightRearStrutPressure"
)
foreach ($name in $strutSensors) {
    $cell = $ws.Columns.Item(1).Find($name)
    if ($cell) {
        $r = $cell.Row
        $ws.Cells.Item($r, 2).Value = 0
        $ws.Cells.Item($r, 3).Value = 30000
    }
}

# --- Update the view to mirror the saved session (selection near the bottom) ---
$lastCell = $ws.Columns.Item(1).Find("EdgeOdometer")
if ($lastCell) {
    $selRow = $lastCell.Row - 1
    $ws.Cells.Item($selRow, 2).Select() | Out-Null
}
